# Update "want-to-go" counts (column F, header "想去人数") on the
# "展览" (Exhibitions), "演出" (Performances) and "全部类型" (All types)
# sheets to reflect a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- 展览 (Exhibitions) sheet ---
$wsExh = $wb.Worksheets.Item("展览")
$wsExh.Range("F2").Value  = 1119
$wsExh.Range("F3").Value  = 249
$wsExh.Range("F4").Value  = 253
$wsExh.Range("F5").Value  = 1793
$wsExh.Range("F6").Value  = 670
$wsExh.Range("F7").Value  = 327
$wsExh.Range("F8").Value  = 520
$wsExh.Range("F9").Value  = 4768
$wsExh.Range("F18").Value = 1824
$wsExh.Range("F19").Value = 112
$wsExh.Range("F22").Value = 52
$wsExh.Range("F24").Value = 940
$wsExh.Range("F25").Value = 308
$wsExh.Range("F27").Value = 2804
$wsExh.Range("F28").Value = 1025
$wsExh.Range("F29").Value = 2488
$wsExh.Range("F31").Value = 1341
$wsExh.Range("F32").Value = 3630
$wsExh.Range("F34").Value = 893
$wsExh.Range("F36").Value = 1138
$wsExh.Range("F37").Value = 932
$wsExh.Range("F38").Value = 1195
$wsExh.Range("F39").Value = 23
$wsExh.Range("F40").Value = 880
$wsExh.Range("F41").Value = 549
$wsExh.Range("F42").Value = 202
$wsExh.Range("F43").Value = 368

# --- 演出 (Performances) sheet ---
$wsPerf = $wb.Worksheets.Item("演出")
$wsPerf.Range("F11").Value = 887

# --- 全部类型 (All types) sheet ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 1119
$wsAll.Range("F3").Value  = 249
$wsAll.Range("F4").Value  = 253
$wsAll.Range("F6").Value  = 1793
$wsAll.Range("F7").Value  = 670
$wsAll.Range("F8").Value  = 327
$wsAll.Range("F9").Value  = 520
$wsAll.Range("F10").Value = 4768
$wsAll.Range("F17").Value = 1824
$wsAll.Range("F18").Value = 112
$wsAll.Range("F21").Value = 887
$wsAll.Range("F24").Value = 52
$wsAll.Range("F26").Value = 940
$wsAll.Range("F27").Value = 308
$wsAll.Range("F28").Value = 2804
$wsAll.Range("F31").Value = 1025
$wsAll.Range("F32").Value = 2488
$wsAll.Range("F33").Value = 1341
$wsAll.Range("F34").Value = 3630
$wsAll.Range("F37").Value = 893
$wsAll.Range("F38").Value = 1138
$wsAll.Range("F39").Value = 932
$wsAll.Range("F41").Value = 1195
$wsAll.Range("F42").Value = 880
$wsAll.Range("F43").Value = 549
$wsAll.Range("F44").Value = 368
